$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.450.46'
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("D3").Value = '1.465.85'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.84%  '
$ws.Range("D5").Value = '0.9160'
$ws.Range("E5").Value = '  -8.57%  '
$ws.Range("D6").Value = '280.55'
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D7").Value = '0.3712'
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").Value = '0.3187'
$ws.Range("E8").Value = '  +3.77%  '
$ws.Range("D9").Value = '40.47'
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D10").Value = '1.053'
$ws.Range("E10").Value = '  +5.27%  '
$ws.Range("D11").Value = '0.06651'
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '5.568'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '18.10'
$ws.Range("E14").Value = '  +6.44%  '
$ws.Range("D15").Value = '6.221'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '1.478.58'
$ws.Range("E16").Value = '  +4.11%  '
$ws.Range("E17").Value = '  +3.04%  '
$ws.Range("D18").Value = '0.9244'
$ws.Range("E18").Value = '  -7.72%  '
$ws.Range("D19").Value = '0.05727'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '71.53'
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("D21").Value = '5.693'
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("D22").Value = '14.71'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("D24").Value = '2.293'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Value = '20.633.17'
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("D26").Value = '2.300'
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("D27").Value = '138.12'
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").Value = '17.52'
$ws.Range("E28").Value = '  +3.46%  '
$ws.Range("D29").Value = '1.635.85'
$ws.Range("E29").Value = '  +3.53%  '
$ws.Range("D30").Value = '113.46'
$ws.Range("E30").Value = '  +4.33%  '
$ws.Range("D31").Value = '3.969'
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("D32").Value = '5.280'
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("D33").Value = '0.8459'
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '0.07803'
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.563'
$ws.Range("E35").Value = '  +13.33%  '
$ws.Range("D36").Value = '0.06094'
$ws.Range("E36").Value = '  +6.27%  '
$ws.Range("D37").Value = '4.881'
$ws.Range("E37").Value = '  +2.43%  '
$ws.Range("D38").Value = '1.154'
$ws.Range("E38").Value = '  +8.68%  '
$ws.Range("D39").Value = '10.66'
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").Value = '0.02064'
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").Value = '0.1893'
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("D42").Value = '0.9414'
$ws.Range("E42").Value = '  -5.99%  '
$ws.Range("D43").Value = '7.373'
$ws.Range("E43").Value = '  -12.28%  '
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("D45").Value = '3.590'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").Value = '12.41'
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").Value = '123.72'
$ws.Range("E47").Value = '  +13.24%  '
$ws.Range("D48").Value = '0.5304'
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("D50").Value = '0.06442'
$ws.Range("E50").Value = '  +4.81%  '
$ws.Range("D51").Value = '1.043'
$ws.Range("E51").Value = '  -0.32%  '
